$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$sheet1 = $wb.Worksheets.Item("Sheet1")
$clientes = $wb.Worksheets.Item("Clientes")

# The buggy flow had created a stray empty "Sheet1" alongside "Clientes".
# Drop the duplicate "Clientes" tab and turn the original first sheet into
# the one-and-only "Clientes" sheet, then rebuild its data with the
# corrected client records.
[void]$clientes.Delete()
$sheet1.Name = "Clientes"
$ws = $sheet1

# Every data cell in this sheet is stored as TEXT (numberStoredAsText),
# so force Text format before writing values — otherwise numeric-looking
# strings (CPF/CEP/phone numbers with leading zeros, etc.) would silently
# get coerced into numbers and lose their leading zeros.
$ws.Range("A1:H3").NumberFormat = "@"

# Header row.
$ws.Range("A1").Value = "Nome"
$ws.Range("B1").Value = "Senha"
$ws.Range("C1").Value = "CPF"
$ws.Range("D1").Value = "Endereço"
$ws.Range("E1").Value = "CEP"
$ws.Range("F1").Value = "Email"
$ws.Range("G1").Value = "Telefone"
$ws.Range("H1").Value = "Endereco"

# Row 2 — fix the bogus placeholder values saved by the buggy "save client" flow.
$ws.Range("A2").Value = "joao"
$ws.Range("B2").Value = "12312312"
$ws.Range("C2").Value = "0000000000"
$ws.Range("E2").Value = "010101010101"
$ws.Range("F2").Value = "brunofraga@gmail.com"
$ws.Range("G2").Value = "3332113221"
$ws.Range("H2").Value = "mato grande"

# Row 3 — same fix.
$ws.Range("A3").Value = "bruno"
$ws.Range("B3").Value = "Itried1993"
$ws.Range("C3").Value = "02370945095"
$ws.Range("E3").Value = "92320-195"
$ws.Range("F3").Value = "brunofraga@gmail.com"
$ws.Range("G3").Value = "51989043802"
$ws.Range("H3").Value = "R. 3 Pinheiros I - Mato Grande"

$excel.DisplayAlerts = $true
